# ---------------------------------------------------------------------------
# Applies the "Employee list is progress" commit to Keywords.xlsx:
#  - LoginPage: header row gets bold/size-12 formatting, selection moves
#  - AddEmployeePage: new "DashboardPage"/"performPIM" column inserted,
#    becomes the active/selected sheet
#  - New EmployeeListPage sheet added with LoginPage+AddEmployeePage-style
#    keyword data plus a DashboardPage column
#  - workbook absPath casing tweak or trivia left to workbook-level props
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$loginPage       = $wb.Worksheets.Item("LoginPage")
$addEmployeePage = $wb.Worksheets.Item("AddEmployeePage")

# --- LoginPage ("sheet1") ---------------------------------------------------
# Fix header casing and apply the bold / size-12 header style (matches the
# style already used on AddEmployeePage).
$loginPage.Range("A1").Value = "TestCaseNo"

$headerRow1 = $loginPage.Range("A1:G1")
$headerRow1.Font.Bold = $true
$headerRow1.Font.Size = 12
$loginPage.Rows.Item(1).RowHeight = 15.6

# Re-fit the columns for the new header font.
$loginPage.Columns.Item(1).ColumnWidth = 11
$loginPage.Columns.Item(2).ColumnWidth = 17.5
$loginPage.Columns.Item(3).ColumnWidth = 16.66667
$loginPage.Columns.Item(4).ColumnWidth = 15.33333
$loginPage.Columns.Item(5).ColumnWidth = 15.33333
$loginPage.Columns.Item(6).ColumnWidth = 12.33333
$loginPage.Columns.Item(7).ColumnWidth = 16.66667

# Move the selection.
$selA12 = $loginPage.Range("A12")
$selA12.Select() | Out-Null

# --- AddEmployeePage ("sheet2") --------------------------------------------
# Insert a new column E carrying the DashboardPage / performPIM keyword data
# (old E:G shift right to F:H).
$addEmployeePage.Columns.Item(5).Insert()

$addEmployeePage.Range("E1").Value = "pages.DashboardPage"
$addEmployeePage.Range("E2").Value = "performPIM"

# Materialise the otherwise-empty E3:E8 cells (present but valueless in the
# target workbook) without pulling in any new style/number-format.
$emptyKeywordRange = $addEmployeePage.Range("E3:E8")
$emptyKeywordRange.Borders.LineStyle = -4142

# Re-fit columns for the extra column.
$addEmployeePage.Columns.Item(1).ColumnWidth = 11
$addEmployeePage.Columns.Item(2).ColumnWidth = 16.66667
$addEmployeePage.Columns.Item(3).ColumnWidth = 16.66667
$addEmployeePage.Columns.Item(4).ColumnWidth = 15.33333
$addEmployeePage.Columns.Item(5).ColumnWidth = 20.66667
$addEmployeePage.Columns.Item(6).ColumnWidth = 23.33333
$addEmployeePage.Columns.Item(7).ColumnWidth = 12.33333
$addEmployeePage.Columns.Item(8).ColumnWidth = 16.66667

# --- EmployeeListPage (new "sheet3") ----------------------------------------
$employeeListPage = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$employeeListPage.Name = "EmployeeListPage"

$employeeListPage.Range("A1").Value = "TestCaseNo"
$employeeListPage.Range("B1").Value = "utils.BrowserUtils"
$employeeListPage.Range("C1").Value = "utils.BrowserUtils"
$employeeListPage.Range("D1").Value = "pages.LoginPage"
$employeeListPage.Range("E1").Value = "pages.DashboardPage"

$employeeListPage.Range("A1:D1").Font.Bold = $true
$employeeListPage.Range("A1:D1").Font.Size = 12
$employeeListPage.Range("E1").Font.Bold = $true
$employeeListPage.Range("E1").Font.Size = 11
$employeeListPage.Rows.Item(1).RowHeight = 15.6

$testCaseNames = @("testCase01","testCase02","testCase03","testCase04","testCase05","testCase06","testCase07")
for ($i = 0; $i -lt $testCaseNames.Count; $i++) {
    $r = $i + 2
    $employeeListPage.Cells.Item($r, 1).Value = $testCaseNames[$i]
    $employeeListPage.Cells.Item($r, 2).Value = "openBrowser"
    $employeeListPage.Cells.Item($r, 3).Value = "launchUrl"
    $employeeListPage.Cells.Item($r, 4).Value = "login"
}
$employeeListPage.Range("E2").Value = "performPIM"

$employeeListPage.Columns.Item(1).ColumnWidth = 11
$employeeListPage.Columns.Item(2).ColumnWidth = 16.66667
$employeeListPage.Columns.Item(3).ColumnWidth = 16.66667
$employeeListPage.Columns.Item(4).ColumnWidth = 15.33333
$employeeListPage.Columns.Item(5).ColumnWidth = 19.77734

$selD2 = $employeeListPage.Range("D2")
$selD2.Select() | Out-Null

# AddEmployeePage is the workbook's active / selected tab in the target
# state, so (re-)activate it last, once every sheet's own selection has
# already been set.
$addEmployeePage.Activate() | Out-Null
$selE2 = $addEmployeePage.Range("E2")
$selE2.Select() | Out-Null

Write-Host "Applied Employee list changes"
